$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Comentarios" column (K) for several requirement rows to "Listo".
$ws.Range("K12").Value = "Listo"
$ws.Range("K13").Value = "Listo"
$ws.Range("K15").Value = "Listo"
$ws.Range("K17").Value = "Listo"
$ws.Range("K18").Value = "Listo"

# Remove the now-unused comment string "Listo, falta recuperar tarea cancelada"
# from the workbook by ensuring no cell references it (it was only used by K12,
# which has already been overwritten above).

# Update the view/selection to match where the user ended up working.
$ws.Application.Goto($ws.Range("A20"), $true)
$ws.Range("K28").Select()
